$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="69.350.32"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  -2.98%  "

$ws.Range("D3").Formula = '="3.683.40"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  -3.41%  "

$ws.Range("D4").Formula = '="1.00"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Formula = '="682.95"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)

$ws.Range("D6").Formula = '="159.99"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -6.73%  "

$ws.Range("D7").Formula = '="3.681.19"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -3.51%  "

$ws.Range("D8").Formula = '="1.00"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -6.03%  "

$ws.Range("E10").Value = "  -9.04%  "

$ws.Range("D11").Formula = '="7.19"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -3.42%  "

$ws.Range("E12").Value = "  -10.33%  "

$ws.Range("E13").Value = "  -6.91%  "

$ws.Range("D14").Formula = '="4.305.41"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -3.20%  "

$ws.Range("D15").Formula = '="32.53"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -10.83%  "

$ws.Range("D16").Formula = '="3.687.23"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Formula = '="69.401.52"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -3.09%  "

$ws.Range("E18").Value = "  -1.24%  "

$ws.Range("D19").Formula = '="15.87"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -9.40%  "

$ws.Range("E20").Value = "  -10.82%  "

$ws.Range("D21").Formula = '="473.92"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -8.07%  "

$ws.Range("E22").Value = "  -5.70%  "

$ws.Range("E23").Value = "  -9.03%  "

$ws.Range("D24").Formula = '="79.46"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -5.47%  "

$ws.Range("D25").Formula = '="3.826.98"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -3.14%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  -11.48%  "

$ws.Range("D28").Formula = '="10.95"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -12.97%  "

$ws.Range("D29").Formula = '="9.21"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -10.73%  "

$ws.Range("D30").Formula = '="2.70"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -10.57%  "

$ws.Range("E31").Value = "  -14.04%  "

$ws.Range("E32").Value = "  -9.81%  "

$ws.Range("E33").Value = "  -10.34%  "

$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("D35").Formula = '="26.70"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -8.61%  "

$ws.Range("E36").Value = "  -7.25%  "

$ws.Range("D37").Formula = '="8.18"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -12.35%  "

$ws.Range("D38").Formula = '="6.09"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -6.69%  "

$ws.Range("D39").Formula = '="2.25"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -7.93%  "

$ws.Range("D41").Formula = '="0.0904"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -10.82%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  -7.05%  "

$ws.Range("D44").Formula = '="165.26"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -1.06%  "

$ws.Range("D45").Formula = '="47.91"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -4.68%  "

$ws.Range("E46").Value = "  -15.95%  "

$ws.Range("E47").Value = "  -6.14%  "

$ws.Range("D48").Formula = '="0.000275"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -10.15%  "

$ws.Range("E49").Value = "  -5.02%  "

$ws.Range("D50").Formula = '="28.30"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -6.50%  "

$ws.Range("E51").Value = "  -9.17%  "

$excel.CutCopyMode = $false
